$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H17").Value = 1667414.9
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 1724901.6
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 5174704.800000001
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -5175040.800000001

$ws.Range("H94").Value = 5697.8667
$ws.Range("I94").Value = 2924
$ws.Range("K94").Value = 2924
$ws.Range("M94").Value = -2473

$ws.Range("H112").Value = 1137.2727
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 1125.3658
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 3376.0974
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -5592.097400000001

$ws.Range("H125").Value = 736659.5600000001
$ws.Range("I125").Value = 1901.7142
$ws.Range("J125").Value = 1308137.9
$ws.Range("K125").Value = 17115.4278
$ws.Range("L125").Value = 11773241.1
$ws.Range("M125").Value = -14655.4278
$ws.Range("N125").Value = -11778161.1

$ws.Range("H137").Value = 1067.7123
$ws.Range("I137").Value = 882.9091
$ws.Range("J137").Value = 1632.3889
$ws.Range("K137").Value = 2648.7273
$ws.Range("L137").Value = 4897.1667
$ws.Range("M137").Value = -98.72730000000001
$ws.Range("N137").Value = -9997.1667

$ws.Range("H138").Value = 2820.7
$ws.Range("I138").Value = 1382.6364
$ws.Range("J138").Value = 5080.514
$ws.Range("K138").Value = 4147.9092
$ws.Range("L138").Value = 15241.542
$ws.Range("M138").Value = 992.0907999999999
$ws.Range("N138").Value = -25521.542

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 1505.05
$ws.Range("J2").Value = 1250
$ws.Range("L2").Value = 1250
$ws.Range("N2").Value = -1476

$ws.Range("H32").Value = 1306.33
$ws.Range("I32").Value = 1165.1075
$ws.Range("K32").Value = 1165.1075
$ws.Range("M32").Value = -878.1075000000001

$ws.Range("H61").Value = 5756
$ws.Range("I61").Value = 8187.5625
$ws.Range("J61").Value = 2219.182
$ws.Range("K61").Value = 8187.5625
$ws.Range("L61").Value = 2219.182
$ws.Range("M61").Value = -7975.5625
$ws.Range("N61").Value = -2643.182

$ws.Range("H92").Value = 19736.125
$ws.Range("J92").Value = 19736.125
$ws.Range("L92").Value = 19736.125
$ws.Range("N92").Value = -24728.125

$ws.Range("H102").Value = 1516
$ws.Range("I102").Value = 1516
$ws.Range("K102").Value = 1516
$ws.Range("M102").Value = 106

$ws.Range("H116").Value = 1505.05
$ws.Range("J116").Value = 1250
$ws.Range("L116").Value = 1250
$ws.Range("N116").Value = -5838

$ws.Range("H136").Value = 5756
$ws.Range("I136").Value = 8187.5625
$ws.Range("J136").Value = 2219.182
$ws.Range("K136").Value = 24562.6875
$ws.Range("L136").Value = 6657.545999999999
$ws.Range("M136").Value = -22012.6875
$ws.Range("N136").Value = -11757.546

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 1505.05
$ws.Range("J3").Value = 1250
$ws.Range("L3").Value = 1250
$ws.Range("N3").Value = -1478

$ws.Range("H20").Value = 7437.8237
$ws.Range("I20").Value = 2032.1428
$ws.Range("J20").Value = 11221.8
$ws.Range("K20").Value = 2032.1428
$ws.Range("L20").Value = 11221.8
$ws.Range("M20").Value = -1785.1428
$ws.Range("N20").Value = -11715.8

$ws.Range("H134").Value = 1796.3062
$ws.Range("I134").Value = 1313.6552
$ws.Range("J134").Value = 2496.15
$ws.Range("K134").Value = 3940.9656
$ws.Range("L134").Value = 7488.450000000001
$ws.Range("M134").Value = -1405.9656
$ws.Range("N134").Value = -12558.45

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 2336.8542
$ws.Range("I31").Value = 1540.3226
$ws.Range("J31").Value = 3789.353
$ws.Range("K31").Value = 1540.3226
$ws.Range("L31").Value = 3789.353
$ws.Range("M31").Value = -1245.3226
$ws.Range("N31").Value = -4379.353

$ws.Range("H34").Value = 2336.8542
$ws.Range("I34").Value = 1540.3226
$ws.Range("J34").Value = 3789.353
$ws.Range("K34").Value = 1540.3226
$ws.Range("L34").Value = 3789.353
$ws.Range("M34").Value = -1338.3226
$ws.Range("N34").Value = -4193.353

$ws = $wb.Worksheets("CUL")
$ws.Range("H14").Value = 312
$ws.Range("I14").Value = 312
$ws.Range("K14").Value = 936
$ws.Range("M14").Value = -763

$ws.Range("H20").Value = 4519.2
$ws.Range("J20").Value = 4519.2
$ws.Range("L20").Value = 13557.6
$ws.Range("N20").Value = -14011.6

$ws.Range("H22").Value = 1382.6086
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 2714.2856
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 8142.8568
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -8480.856800000001

$ws.Range("H27").Value = 1382.6086
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 2714.2856
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 8142.8568
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -8346.856800000001

$ws.Range("H32").Value = 500000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H34").Value = 599
$ws.Range("I34").Value = 598
$ws.Range("K34").Value = 1794
$ws.Range("M34").Value = -1710

$ws.Range("H39").Value = 9585.296
$ws.Range("J39").Value = 6192
$ws.Range("L39").Value = 18576
$ws.Range("N39").Value = -19164

$ws.Range("H46").Value = 110
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 600
$ws.Range("M46").Value = -209
$ws.Range("N46").Value = -782

$ws.Range("H55").Value = 12465
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 12465
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 37395
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -37749

$ws.Range("H58").Value = 1166.6666
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 2250
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -2122
$ws.Range("N58").Value = -6256

$ws.Range("H131").Value = 5314.5454
$ws.Range("J131").Value = 5926.207
$ws.Range("L131").Value = 17778.621
$ws.Range("N131").Value = -27858.621

$ws.Range("H140").Value = 2267.3
$ws.Range("I140").Value = 1902.1666
$ws.Range("J140").Value = 2815
$ws.Range("K140").Value = 5706.4998
$ws.Range("L140").Value = 8445
$ws.Range("M140").Value = -526.4997999999996
$ws.Range("N140").Value = -18805

$ws = $wb.Worksheets("GSM")
$ws.Range("H113").Value = 1581.875
$ws.Range("I113").Value = 1592.5
$ws.Range("J113").Value = 1550
$ws.Range("K113").Value = 1592.5
$ws.Range("L113").Value = 1550
$ws.Range("M113").Value = 577.5
$ws.Range("N113").Value = -5890

$ws = $wb.Worksheets("LTW")
$ws.Range("H46").Value = 250750.25
$ws.Range("I46").Value = 334000.34
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 334000.34
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -333812.34
$ws.Range("N46").Value = -1376

$ws = $wb.Worksheets("WVR")
$ws.Range("H96").Value = 1650
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1650
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1650
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4396

$ws.Range("H107").Value = 2895.8333
$ws.Range("I107").Value = 3250
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 9750
$ws.Range("L107").Value = 7200
$ws.Range("M107").Value = -7830
$ws.Range("N107").Value = -11040

$ws.Range("H113").Value = 2820
$ws.Range("I113").Value = 3025
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 9075
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -6905
$ws.Range("N113").Value = -10340

$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

$ws.Range("H136").Value = 3898.7144
$ws.Range("I136").Value = 583.4286
$ws.Range("J136").Value = 10529.286
$ws.Range("K136").Value = 1750.2858
$ws.Range("L136").Value = 31587.858
$ws.Range("M136").Value = 799.7142000000001
$ws.Range("N136").Value = -36687.858
